$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.557.90"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "3.015.01"
$ws.Range("E3").Value = "  +1.94%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'378.61"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").Value = "'102.28"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("D10").Value = "'36.66"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "3.495.40"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").Value = "'18.38"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "3.024.45"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").Value = "'0.971"
$ws.Range("E17").Value = "  -4.14%  "

$ws.Range("D18").Value = "'10.61"
$ws.Range("E18").Value = "  -14.93%  "

$ws.Range("D19").Value = "51.546.67"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Value = "'3.08"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "'12.40"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("D23").Value = "'69.87"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").Value = "'266.16"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  -7.48%  "

$ws.Range("E26").Value = "  +3.37%  "

$ws.Range("D27").Value = "'7.57"
$ws.Range("E27").Value = "  +8.77%  "

$ws.Range("E28").Value = "  +4.18%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'26.13"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").Value = "'0.107"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").Value = "'10.23"
$ws.Range("E32").Value = "  -2.64%  "

$ws.Range("D33").Value = "'2.05"
$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").Value = "'50.53"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").Value = "'33.75"
$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").Value = "'0.0449"
$ws.Range("E36").Value = "  +3.14%  "

$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  +2.46%  "

$ws.Range("D39").Value = "'0.288"
$ws.Range("E39").Value = "  +11.52%  "

$ws.Range("D40").Value = "'16.89"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("E41").Value = "  +1.40%  "

$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.53"
$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.76"
$ws.Range("E44").Value = "  +5.62%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'123.74"
$ws.Range("E45").Value = "  +2.75%  "

$ws.Range("D46").Value = "'21.60"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'2.09"
$ws.Range("E47").Value = "  +3.00%  "

$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("D49").Value = "2.024.81"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").Value = "3.316.23"
$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("D51").Value = "'0.0317"
$ws.Range("E51").Value = "  -1.25%  "

